$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2880.7568
$ws.Range("J112").Value = 3046.7058
$ws.Range("L112").Value = 9140.117400000001
$ws.Range("N112").Value = -11356.1174

$ws.Range("H129").Value = 1656.3043
$ws.Range("I129").Value = 1098.2
$ws.Range("J129").Value = 1811.3334
$ws.Range("K129").Value = 3294.6
$ws.Range("L129").Value = 5434.0002
$ws.Range("M129").Value = 1705.4
$ws.Range("N129").Value = -15434.0002

$ws.Range("H132").Value = 5942.737
$ws.Range("I132").Value = 5751.1763
$ws.Range("J132").Value = 6225.913
$ws.Range("K132").Value = 17253.5289
$ws.Range("L132").Value = 18677.739
$ws.Range("M132").Value = -14723.5289
$ws.Range("N132").Value = -23737.739

$ws.Range("H137").Value = 1365.1864
$ws.Range("I137").Value = 1673.2188
$ws.Range("K137").Value = 5019.6564
$ws.Range("M137").Value = -2469.6564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9883.763000000001
$ws.Range("I32").Value = 5461.988
$ws.Range("J32").Value = 34056.133
$ws.Range("K32").Value = 5461.988
$ws.Range("L32").Value = 34056.133
$ws.Range("M32").Value = -5174.988
$ws.Range("N32").Value = -34630.133

$ws.Range("H61").Value = 2536.9363
$ws.Range("I61").Value = 3547.0833
$ws.Range("J61").Value = 1482.8695
$ws.Range("K61").Value = 3547.0833
$ws.Range("L61").Value = 1482.8695
$ws.Range("M61").Value = -3335.0833
$ws.Range("N61").Value = -1906.8695

$ws.Range("H74").Value = 1441.4833
$ws.Range("I74").Value = 1033.925
$ws.Range("J74").Value = 2256.6
$ws.Range("K74").Value = 1033.925
$ws.Range("L74").Value = 2256.6
$ws.Range("M74").Value = -159.925
$ws.Range("N74").Value = -4004.6

$ws.Range("H77").Value = 1441.4833
$ws.Range("I77").Value = 1033.925
$ws.Range("J77").Value = 2256.6
$ws.Range("K77").Value = 5169.625
$ws.Range("L77").Value = 11283
$ws.Range("M77").Value = -801.625
$ws.Range("N77").Value = -20019

$ws.Range("H88").Value = 1776
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 2490
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 2490
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -3302

$ws.Range("H91").Value = 1776
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 2490
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 2490
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -5298

$ws.Range("H122").Value = 1511.1818
$ws.Range("I122").Value = 1128.4
$ws.Range("J122").Value = 2331.4285
$ws.Range("K122").Value = 3385.2
$ws.Range("L122").Value = 6994.2855
$ws.Range("M122").Value = -935.2000000000003
$ws.Range("N122").Value = -11894.2855

$ws.Range("H136").Value = 2536.9363
$ws.Range("I136").Value = 3547.0833
$ws.Range("J136").Value = 1482.8695
$ws.Range("K136").Value = 10641.2499
$ws.Range("L136").Value = 4448.6085
$ws.Range("M136").Value = -8091.249899999999
$ws.Range("N136").Value = -9548.6085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1648.4546
$ws.Range("I86").Value = 1242.1666
$ws.Range("J86").Value = 2136
$ws.Range("K86").Value = 1242.1666
$ws.Range("L86").Value = 2136
$ws.Range("M86").Value = -119.1666
$ws.Range("N86").Value = -4382

$ws.Range("H89").Value = 1648.4546
$ws.Range("I89").Value = 1242.1666
$ws.Range("J89").Value = 2136
$ws.Range("K89").Value = 6210.833000000001
$ws.Range("L89").Value = 10680
$ws.Range("M89").Value = -594.8330000000005
$ws.Range("N89").Value = -21912

$ws.Range("H134").Value = 3663.7715
$ws.Range("I134").Value = 1620.8684
$ws.Range("J134").Value = 6089.7188
$ws.Range("K134").Value = 4862.6052
$ws.Range("L134").Value = 18269.1564
$ws.Range("M134").Value = -2327.6052
$ws.Range("N134").Value = -23339.1564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7754163
$ws.Range("I31").Value = 1875.0344
$ws.Range("K31").Value = 1875.0344
$ws.Range("M31").Value = -1580.0344

$ws.Range("H34").Value = 7754163
$ws.Range("I34").Value = 1875.0344
$ws.Range("K34").Value = 1875.0344
$ws.Range("M34").Value = -1673.0344

$ws.Range("H58").Value = 4202.1914
$ws.Range("I58").Value = 1993.2759
$ws.Range("J58").Value = 7761
$ws.Range("K58").Value = 1993.2759
$ws.Range("L58").Value = 7761
$ws.Range("M58").Value = -1790.2759
$ws.Range("N58").Value = -8167

$ws.Range("H62").Value = 18808.75
$ws.Range("I62").Value = 13378.333
$ws.Range("J62").Value = 35100
$ws.Range("K62").Value = 13378.333
$ws.Range("L62").Value = 35100
$ws.Range("M62").Value = -12754.333
$ws.Range("N62").Value = -36348

$ws.Range("H65").Value = 18808.75
$ws.Range("I65").Value = 13378.333
$ws.Range("J65").Value = 35100
$ws.Range("K65").Value = 66891.66500000001
$ws.Range("L65").Value = 175500
$ws.Range("M65").Value = -63771.66500000001
$ws.Range("N65").Value = -181740

$ws.Range("H132").Value = 1818.4921
$ws.Range("I132").Value = 1072.7142
$ws.Range("K132").Value = 3218.1426
$ws.Range("M132").Value = -688.1425999999997

$ws.Range("H134").Value = 1677.6531
$ws.Range("I134").Value = 915.7083
$ws.Range("J134").Value = 2409.12
$ws.Range("K134").Value = 2747.1249
$ws.Range("L134").Value = 7227.36
$ws.Range("M134").Value = -212.1248999999998
$ws.Range("N134").Value = -12297.36

$ws.Range("H136").Value = 4202.1914
$ws.Range("I136").Value = 1993.2759
$ws.Range("J136").Value = 7761
$ws.Range("K136").Value = 5979.8277
$ws.Range("L136").Value = 23283
$ws.Range("M136").Value = -3429.8277
$ws.Range("N136").Value = -28383

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8200
$ws.Range("I70").Value = 7500
$ws.Range("J70").Value = 8666.666999999999
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 8666.666999999999
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -9206.666999999999

$ws.Range("H73").Value = 8200
$ws.Range("I73").Value = 7500
$ws.Range("J73").Value = 8666.666999999999
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 8666.666999999999
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -10538.667

$ws.Range("H92").Value = 3845.9167
$ws.Range("J92").Value = 3845.9167
$ws.Range("L92").Value = 3845.9167
$ws.Range("N92").Value = -7589.9167

$ws.Range("H122").Value = 3128.7083
$ws.Range("I122").Value = 3056.05
$ws.Range("J122").Value = 3492
$ws.Range("K122").Value = 9168.150000000001
$ws.Range("L122").Value = 10476
$ws.Range("M122").Value = -6718.150000000001
$ws.Range("N122").Value = -15376

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16968676
$ws.Range("I132").Value = 28602936
$ws.Range("J132").Value = 2046
$ws.Range("K132").Value = 85808808
$ws.Range("L132").Value = 6138
$ws.Range("M132").Value = -85806278
$ws.Range("N132").Value = -11198

$ws.Range("H136").Value = 10871149
$ws.Range("I136").Value = 15625967
$ws.Range("J136").Value = 2994.2856
$ws.Range("K136").Value = 46877901
$ws.Range("L136").Value = 8982.856800000001
$ws.Range("M136").Value = -46875351
$ws.Range("N136").Value = -14082.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 973.7857
$ws.Range("I113").Value = 475.5
$ws.Range("J113").Value = 1347.5
$ws.Range("K113").Value = 1426.5
$ws.Range("L113").Value = 4042.5
$ws.Range("M113").Value = 743.5
$ws.Range("N113").Value = -8382.5

$ws.Range("H132").Value = 1976.8214
$ws.Range("I132").Value = 1866.7667
$ws.Range("J132").Value = 2103.8076
$ws.Range("K132").Value = 5600.300099999999
$ws.Range("L132").Value = 6311.4228
$ws.Range("M132").Value = -3070.300099999999
$ws.Range("N132").Value = -11371.4228

$ws.Range("H136").Value = 4908027.5
$ws.Range("I136").Value = 7150620
$ws.Range("J136").Value = 2356.25
$ws.Range("K136").Value = 21451860
$ws.Range("L136").Value = 7068.75
$ws.Range("M136").Value = -21449310
$ws.Range("N136").Value = -12168.75
